# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, column letter, new text value.
$updates = @(
    @{Row=2; Col='D'; Value='26.805.65'},
    @{Row=2; Col='E'; Value='  +0.10%  '},
    @{Row=3; Col='D'; Value='1.649.55'},
    @{Row=3; Col='E'; Value='  -0.04%  '},
    @{Row=5; Col='D'; Value='216.83'},
    @{Row=5; Col='E'; Value='  +0.71%  '},
    @{Row=6; Col='D'; Value='0.506'},
    @{Row=6; Col='E'; Value='  +0.21%  '},
    @{Row=7; Col='E'; Value='  +0.63%  '},
    @{Row=8; Col='D'; Value='0.253'},
    @{Row=8; Col='E'; Value='  +0.19%  '},
    @{Row=9; Col='D'; Value='0.0628'},
    @{Row=9; Col='E'; Value='  -0.26%  '},
    @{Row=10; Col='E'; Value='  -0.07%  '},
    @{Row=11; Col='E'; Value='  +0.20%  '},
    @{Row=12; Col='D'; Value='1.873.53'},
    @{Row=12; Col='E'; Value='  -0.32%  '},
    @{Row=13; Col='D'; Value='1.657.46'},
    @{Row=13; Col='E'; Value='  +0.12%  '},
    @{Row=14; Col='E'; Value='  +1.50%  '},
    @{Row=15; Col='E'; Value='  +0.17%  '},
    @{Row=16; Col='D'; Value='65.77'},
    @{Row=16; Col='E'; Value='  -0.71%  '},
    @{Row=17; Col='D'; Value='26.807.50'},
    @{Row=17; Col='E'; Value='  -0.07%  '},
    @{Row=18; Col='E'; Value='  -0.24%  '},
    @{Row=19; Col='D'; Value='216.96'},
    @{Row=19; Col='E'; Value='  -0.85%  '},
    @{Row=20; Col='E'; Value='  +0.71%  '},
    @{Row=21; Col='E'; Value='  +0.23%  '},
    @{Row=22; Col='E'; Value='  +16.00%  '},
    @{Row=23; Col='E'; Value='  -0.89%  '},
    @{Row=24; Col='D'; Value='9.48'},
    @{Row=24; Col='E'; Value='  -0.04%  '},
    @{Row=25; Col='D'; Value='145.76'},
    @{Row=25; Col='E'; Value='  -1.34%  '},
    @{Row=26; Col='E'; Value='  +0.75%  '},
    @{Row=27; Col='D'; Value='0.121'},
    @{Row=27; Col='E'; Value='  -0.66%  '},
    @{Row=28; Col='E'; Value='  +3.75%  '},
    @{Row=29; Col='E'; Value='  +0.04%  '},
    @{Row=30; Col='E'; Value='  -0.05%  '},
    @{Row=31; Col='E'; Value='  +0.62%  '},
    @{Row=32; Col='E'; Value='  -0.94%  '},
    @{Row=33; Col='E'; Value='  +0.59%  '},
    @{Row=34; Col='E'; Value='  +0.81%  '},
    @{Row=35; Col='D'; Value='1.278.52'},
    @{Row=35; Col='E'; Value='  -0.32%  '},
    @{Row=36; Col='E'; Value='  +1.94%  '},
    @{Row=37; Col='E'; Value='  -0.05%  '},
    @{Row=38; Col='E'; Value='  +4.90%  '},
    @{Row=39; Col='D'; Value='0.834'},
    @{Row=39; Col='E'; Value='  +2.92%  '},
    @{Row=40; Col='E'; Value='  +0.70%  '},
    @{Row=41; Col='D'; Value='0.820'},
    @{Row=41; Col='E'; Value='  +1.39%  '},
    @{Row=42; Col='E'; Value='  -1.56%  '},
    @{Row=43; Col='E'; Value='  +0.85%  '},
    @{Row=44; Col='D'; Value='1.798.98'},
    @{Row=44; Col='E'; Value='  +0.58%  '},
    @{Row=45; Col='D'; Value='92.17'},
    @{Row=45; Col='E'; Value='  -1.68%  '},
    @{Row=46; Col='D'; Value='59.65'},
    @{Row=46; Col='E'; Value='  +6.19%  '},
    @{Row=47; Col='E'; Value='  +1.38%  '},
    @{Row=48; Col='E'; Value='  +1.28%  '},
    @{Row=49; Col='E'; Value='  +0.39%  '},
    @{Row=50; Col='D'; Value='7.78'},
    @{Row=50; Col='E'; Value='  +1.00%  '},
    @{Row=51; Col='D'; Value='0.0985'},
    @{Row=51; Col='E'; Value='  +1.57%  '}
)

foreach ($u in $updates) {
    $cell = $ws.Range("$($u.Col)$($u.Row)")
    if ($u.Col -eq 'D') {
        # Price column holds numeric-looking text (e.g. "216.83", "1.649.55").
        # Force Text formatting first so Excel doesn't coerce it into a
        # Number/Date, then restore the default (unstyled) look.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
